$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 0.6666666666666666
$ws.Range("G2").Value = 0.01427633333333333
$ws.Range("H2").Value = 0.042829
$ws.Range("I2").Value = 0.009697706618844965
$ws.Range("J2").Value = 0.009697706618844965
$ws.Range("M2").Value = 0.06687733333333333
$ws.Range("N2").Value = 0.200632
$ws.Range("O2").Value = 0.2204558290608057
$ws.Range("P2").Value = 0.2204558290608058
$ws.Range("Q2").Value = 0.0009547631031111111
$ws.Range("R2").Value = 0.008592867928
$ws.Range("S2").Value = 0.00213791595264593
$ws.Range("T2").Value = 0.00213791595264593
$ws.Range("E3").Value = 2
$ws.Range("F3").Value = 0.6666666666666666
$ws.Range("G3").Value = 0.01427633333333333
$ws.Range("H3").Value = 0.042829
$ws.Range("I3").Value = 0.009697706618844965
$ws.Range("J3").Value = 0.009697706618844965
$ws.Range("O3").Value = 0.6972292484819982
$ws.Range("P3").Value = 0.6972292484819982
$ws.Range("Q3").Value = 0.003019601539666667
$ws.Range("R3").Value = 0.027176413857
$ws.Range("S3").Value = 0.006761524697856174
$ws.Range("T3").Value = 0.006761524697856174
$ws.Range("E4").Value = 2
$ws.Range("F4").Value = 0.6666666666666666
$ws.Range("G4").Value = 0.01427633333333333
$ws.Range("H4").Value = 0.042829
$ws.Range("I4").Value = 0.009697706618844965
$ws.Range("J4").Value = 0.009697706618844965
$ws.Range("O4").Value = 0.08231492245719596
$ws.Range("P4").Value = 0.08231492245719596
$ws.Range("Q4").Value = 0.0003564943196666666
$ws.Range("R4").Value = 0.003208448877
$ws.Range("S4").Value = 0.0007982659683428593
$ws.Range("T4").Value = 0.0007982659683428593
$ws.Range("I5").Value = 0.7994031344498523
$ws.Range("J5").Value = 0.7994031344498523
$ws.Range("M5").Value = 0.06687733333333333
$ws.Range("N5").Value = 0.200632
$ws.Range("O5").Value = 0.2204558290608057
$ws.Range("P5").Value = 0.2204558290608058
$ws.Range("Q5").Value = 0.07870320760177778
$ws.Range("R5").Value = 0.7083288684160001
$ws.Range("S5").Value = 0.1762330807589489
$ws.Range("T5").Value = 0.176233080758949
$ws.Range("I6").Value = 0.7994031344498523
$ws.Range("J6").Value = 0.7994031344498523
$ws.Range("O6").Value = 0.6972292484819982
$ws.Range("P6").Value = 0.6972292484819982
$ws.Range("S6").Value = 0.5573672466666243
$ws.Range("T6").Value = 0.5573672466666243
$ws.Range("I7").Value = 0.7994031344498523
$ws.Range("J7").Value = 0.7994031344498523
$ws.Range("O7").Value = 0.08231492245719596
$ws.Range("P7").Value = 0.08231492245719596
$ws.Range("S7").Value = 0.06580280702427899
$ws.Range("T7").Value = 0.06580280702427899
$ws.Range("G8").Value = 0.2810293333333334
$ws.Range("H8").Value = 0.8430880000000001
$ws.Range("I8").Value = 0.1908991589313027
$ws.Range("J8").Value = 0.1908991589313027
$ws.Range("M8").Value = 0.06687733333333333
$ws.Range("N8").Value = 0.200632
$ws.Range("O8").Value = 0.2204558290608057
$ws.Range("P8").Value = 0.2204558290608058
$ws.Range("Q8").Value = 0.01879449240177778
$ws.Range("R8").Value = 0.169150431616
$ws.Range("S8").Value = 0.04208483234921086
$ws.Range("T8").Value = 0.04208483234921087
$ws.Range("G9").Value = 0.2810293333333334
$ws.Range("H9").Value = 0.8430880000000001
$ws.Range("I9").Value = 0.1908991589313027
$ws.Range("J9").Value = 0.1908991589313027
$ws.Range("O9").Value = 0.6972292484819982
$ws.Range("P9").Value = 0.6972292484819982
$ws.Range("Q9").Value = 0.05944079532266667
$ws.Range("R9").Value = 0.5349671579040001
$ws.Range("S9").Value = 0.1331004771175177
$ws.Range("T9").Value = 0.1331004771175177
$ws.Range("G10").Value = 0.2810293333333334
$ws.Range("H10").Value = 0.8430880000000001
$ws.Range("I10").Value = 0.1908991589313027
$ws.Range("J10").Value = 0.1908991589313027
$ws.Range("O10").Value = 0.08231492245719596
$ws.Range("P10").Value = 0.08231492245719596
$ws.Range("Q10").Value = 0.007017583482666666
$ws.Range("S10").Value = 0.01571384946457411
$ws.Range("T10").Value = 0.01571384946457411
